{"js": "// Apply the \"Added many more features\" edits to the Ice Ice Yeti review.\n// Each change is a straightforward text replacement; several target\n// strings are unique, but the title/H1 line (\"Play Ice Ice Yeti for\n// Free - Review of the Arctic-themed Slot Game\") appears twice (once as\n// the Heading1 and once as a bold run near the end) and both need the\n// same replacement, so we replace every match search() returns.\n\nconst replacements = [\n  [\n    \"Play Ice Ice Yeti for Free - Review of the Arctic-themed Slot Game\",\n    \"Play Ice Ice Yeti for Free\"\n  ],\n  [\n    \"Expanding reels offer more gameplay area\",\n    \"Expanding reels for increased gameplay area\"\n  ],\n  [\n    \"Maximum win of 5,000x the bet\",\n    \"Chance to win a maximum of 5,000x the bet\"\n  ],\n  [\n    \"Available on desktop and mobile\",\n    \"Playable on desktop, mobile, and tablet\"\n  ],\n  [\n    \"Jackpot may not be as high as some other games\",\n    \"Limited bonus features\"\n  ],\n  [\n    \"No free spins bonus round\",\n    \"Random triggering of Yeti Shake function\"\n  ],\n  [\n    \"Read our review of Ice Ice Yeti, an Arctic-themed slot game with expanding reels and up to 16,807 ways to win. Play for free on desktop and mobile.\",\n    \"Read our review of Ice Ice Yeti, a slot game with expanding reels and 16,807 ways to win. Play for free.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Ice Ice Yeti review.\n# Each change is a straightforward Find/Replace; the title/H1 line\n# (\"Play Ice Ice Yeti for Free - Review of the Arctic-themed Slot Game\")\n# appears twice (Heading1 + a bold run near the end) and both need the\n# same replacement, so wdReplaceAll is used for every pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Ice Ice Yeti for Free - Review of the Arctic-themed Slot Game\", \"Play Ice Ice Yeti for Free\"),\n    @(\"Expanding reels offer more gameplay area\", \"Expanding reels for increased gameplay area\"),\n    @(\"Maximum win of 5,000x the bet\", \"Chance to win a maximum of 5,000x the bet\"),\n    @(\"Available on desktop and mobile\", \"Playable on desktop, mobile, and tablet\"),\n    @(\"Jackpot may not be as high as some other games\", \"Limited bonus features\"),\n    @(\"No free spins bonus round\", \"Random triggering of Yeti Shake function\"),\n    @(\"Read our review of Ice Ice Yeti, an Arctic-themed slot game with expanding reels and up to 16,807 ways to win. Play for free on desktop and mobile.\", \"Read our review of Ice Ice Yeti, a slot game with expanding reels and 16,807 ways to win. Play for free.\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
